$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 currently holds "2020-12-31" / 0. It becomes the final data row,
# representing "2022-12-31" / 1362 (the rows for 2020 and 2021 are removed).
# Use a formula producing a text string so Excel's autodetection does not
# convert the "YYYY-MM-DD" text into a date serial number, then convert
# that formula result into a plain cached value via copy / paste-values so
# the cell keeps its original style (border/bold/center) untouched.
$ws.Range("A5").Formula = "=""2022-12-31"""
$ws.Range("A5").Copy()
$ws.Range("A5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("B5").Value = 1362

# Remove the now-obsolete rows for 2020-12-31 and 2021-12-31.
$ws.Rows("6:7").Delete()
